$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.356.73"
$ws.Range("E2").Value = "  +2.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.108.20"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "345.43"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5223"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4441"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.39"
$ws.Range("E9").Value = "  +4.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09473"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.10"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.727"
$ws.Range("E13").Value = "  +7.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.136.17"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.911"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "101.70"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001164"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.30"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06734"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.321"
$ws.Range("E21").Value = "  +2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.389.64"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.63"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.371.95"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("E27").Value = "  +0.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.543"
$ws.Range("E28").Value = "  +0.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.29"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.42"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.760"
$ws.Range("E32").Value = "  +8.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1055"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.911"
$ws.Range("E34").Value = "  +13.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.262"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.932"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.53"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02641"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06803"
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7037"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.347"
$ws.Range("E41").Value = "  +4.60%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.56"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2228"
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6841"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.54"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.364"
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.355"
$ws.Range("E48").Value = "  +16.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.648"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000346"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("E51").Value = "  +0.34%  "
